$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Row 2: existing mission record gets a new employee + recalculated expenses
# ---------------------------------------------------------------------------
$ws.Range("B2").Value = 20
$ws.Range("C2").Value = "عبد العزيز صلاح عبد العزيز على حسن"

$expenses2 = '[{"id":"expense-1","type":"transportation","amount":100,"banks":["كريدى","مانى فيللوز"]},{"id":"expense-2","type":"transportation","amount":30,"banks":["اسكندرية"]},{"id":"expense-3","type":"hospitality","amount":0,"banks":[]},{"id":"expense-4","type":"fees","amount":0,"banks":[]}]'
$ws.Range("H2").Value = $expenses2

# I2 holds a number-looking value but must stay text (matches numberStoredAsText
# ignoredError already on the sheet) - force text format, write it, then drop
# the leftover number-format style so the cell keeps default styling.
$ws.Range("I2").NumberFormat = "@"
$ws.Range("I2").Value = "130"
$ws.Range("I2").ClearFormats()

# ---------------------------------------------------------------------------
# Row 3: brand-new mission row (this used to be the only data row's content)
# ---------------------------------------------------------------------------
$ws.Range("A3").Value = "0e3f0a1e-d8a6-44df-b229-03628e5bca16"
$ws.Range("B3").Value = 62
$ws.Range("C3").Value = "محمد مجدى السيد عبد الدايم"
$ws.Range("D3").Value = "20أ القاهرة"

$ws.Range("E3").NumberFormat = "@"
$ws.Range("E3").Value = "2025-09-14"
$ws.Range("E3").ClearFormats()

$expenses3 = '[{"id":"expense-1","type":"fees","amount":50,"banks":["كريدى","مانى فيللوز"]},{"id":"expense-2","type":"fees","amount":50,"banks":["كريدى"]}]'
$ws.Range("H3").Value = $expenses3

$ws.Range("I3").NumberFormat = "@"
$ws.Range("I3").Value = "100"
$ws.Range("I3").ClearFormats()

$ws.Range("J3").Value = "2025-09-14T11:31:03.321Z"

# ---------------------------------------------------------------------------
# Row 4: brand-new, essentially-blank mission row
# ---------------------------------------------------------------------------
$ws.Range("A4").Value = "6c41556f-a921-4153-bbf5-ec5a27eb2838"
$ws.Range("B4").Value = 0

# employeeName / employeeBranch are blank for this placeholder row. A literal
# "" assignment removes the cell entirely in this engine (no way to keep an
# explicit empty text cell without also leaving formula/calcChain residue
# that wouldn't otherwise belong in this workbook), so leave them unset -
# reads back as "" either way.
$ws.Range("C4").Value = ""
$ws.Range("D4").Value = ""

$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "2025-09-14"
$ws.Range("E4").ClearFormats()

$ws.Range("H4").Value = "[]"

$ws.Range("I4").NumberFormat = "@"
$ws.Range("I4").Value = "0"
$ws.Range("I4").ClearFormats()

$ws.Range("J4").Value = "2025-09-14T11:31:55.900Z"
